# Updated cryptos list on Tue May 21 06:31:10 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '71.191.44'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +6.22%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.674.90'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +18.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.57'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +4.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '182.49'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +2.53%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.670.19'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +18.05%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.538'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +4.83%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.163'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +7.68%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.65'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +4.41%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.499'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +5.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.73'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +12.35%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000255'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +5.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.285.25'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +18.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '71.167.00'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +6.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.655.37'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  +17.58%  '
$ws.Range("E18").Value = '  +1.17%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.52'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +7.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.08'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.88%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '520.11'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +6.91%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.26'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +19.14%  '
$ws.Range("E23").Value = '  +7.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.32'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +5.46%  '
$ws.Range("E25").Value = '  +11.00%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.60'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +7.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.05'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +7.27%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.56'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +11.79%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '8.15'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.96%  '
$ws.Range("E31").Value = '  +7.87%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '31.80'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  +12.90%  '
$ws.Range("E33").Value = '  +17.89%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.117'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +3.97%  '
$ws.Range("E35").Value = '  -0.11%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.16'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +9.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.02'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +7.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.347'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +11.86%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.19'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +9.63%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.130'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +5.46%  '
$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '51.15'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +3.90%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '45.13'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -7.63%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.87'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +6.47%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.137.19'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +12.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '417.35'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +11.90%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.81'
$ws.Range("D46").ClearFormats()
$ws.Range("E47").Value = '  +6.89%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '28.44'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +13.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '139.54'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +3.10%  '
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.48'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +10.50%  '
